# Apply "setting taps from banked or tanked regulators" edit.
#
# 1) Transformer sheet: set Tap 1 / Lowest Tap / Highest Tap / Min Range (%) /
#    Max Range (%) for the three single-phase voltage regulators (rows 23-25).
# 2) Bus sheet: the phase angle bookkeeping for several buses had the A/B
#    phase rows swapped (labels + angle values), and the 692 bus had its
#    A/B/C rows rotated by one position. Fix the labels/angles accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Transformer sheet — regulator taps
# ---------------------------------------------------------------------------
$wsXfmr = $wb.Worksheets.Item("Transformer")

# Row 23 -> Reg:reg1
$wsXfmr.Range("P23").Value = 10
$wsXfmr.Range("S23").Value = -16
$wsXfmr.Range("T23").Value = 16
$wsXfmr.Range("U23").Value = 10
$wsXfmr.Range("V23").Value = 10

# Row 24 -> Reg:reg2
$wsXfmr.Range("P24").Value = 8
$wsXfmr.Range("S24").Value = -16
$wsXfmr.Range("T24").Value = 16
$wsXfmr.Range("U24").Value = 10
$wsXfmr.Range("V24").Value = 10

# Row 25 -> Reg:reg3
$wsXfmr.Range("P25").Value = 11
$wsXfmr.Range("S25").Value = -16
$wsXfmr.Range("T25").Value = 16
$wsXfmr.Range("U25").Value = 10
$wsXfmr.Range("V25").Value = 10

# ---------------------------------------------------------------------------
# Bus sheet — fix swapped phase-A / phase-B rows (and the 692 rotation)
# ---------------------------------------------------------------------------
$wsBus = $wb.Worksheets.Item("Bus")

function Set-BusRow($row, $name, $angle) {
    $wsBus.Range("A$row").Value = $name
    $wsBus.Range("E$row").Value = $angle
}

# Simple A/B swaps: (row with former "_A", row with former "_B")
Set-BusRow 4  "632_B"       -120
Set-BusRow 5  "632_A"       0

Set-BusRow 7  "633_B"       -120
Set-BusRow 8  "633_A"       0

Set-BusRow 10 "634_B"       -120
Set-BusRow 11 "634_A"       0

Set-BusRow 17 "650_B"       -120
Set-BusRow 18 "650_A"       0

Set-BusRow 21 "670_B"       -120
Set-BusRow 22 "670_A"       0

Set-BusRow 24 "671_B"       -120
Set-BusRow 25 "671_A"       0

Set-BusRow 27 "675_B"       -120
Set-BusRow 28 "675_A"       0

Set-BusRow 30 "680_B"       -120
Set-BusRow 31 "680_A"       0

# 692_C / 692_A / 692_B rotate down by one row
Set-BusRow 34 "692_B"       -120
Set-BusRow 35 "692_C"       120
Set-BusRow 36 "692_A"       0

Set-BusRow 38 "brkr_B"      -120
Set-BusRow 39 "brkr_A"      0

Set-BusRow 43 "mid_B"       -120
Set-BusRow 44 "mid_A"       0

Set-BusRow 46 "rg60_B"      -120
Set-BusRow 47 "rg60_A"      0

Set-BusRow 49 "sourcebus_B" -120
Set-BusRow 50 "sourcebus_A" 0

Set-BusRow 53 "xf1_B"       -120
Set-BusRow 54 "xf1_A"       0

Write-Output "edit applied"
